# Intervention content and delivery.xlsx
#
# Several rows (entities that have become obsolete / superseded) are
# re-styled to match the workbook's "obsolete" row formatting (the dark
# fill already used elsewhere in the sheet, cellXfs style index 2) and
# their "Curation status" column (S) is updated to "Obsolete".
#
# Row 34 already carries that target style across A:V, so we copy its
# formatting (format-only paste) onto each affected row instead of
# poking Interior.Color directly - that keeps the workbook's existing
# shared style entry instead of Excel minting a brand new (duplicate)
# style/fill for every touched range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$templateRange = "A34:V34"

$rows = @(28, 31, 32, 33, 41, 42, 43, 47, 48, 49, 52, 53, 54, 55, 58, 60, 61, 62, 63)

foreach ($r in $rows) {
    $ws.Range($templateRange).Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)
    $ws.Range("S" + $r).Value = "Obsolete"
}
